# Update cryptocurrency prices in column D (Price) on the active worksheet.
# Values are stored as text (inline strings) in the source workbook, so we
# prefix each value with a leading apostrophe to force Excel to keep the
# cell as text (preserving exact formatting, leading/trailing zeros, etc.)
# instead of coercing it into a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'263.67"
$ws.Range("D3").Value  = "'21.27"
$ws.Range("D4").Value  = "'6.131"
$ws.Range("D5").Value  = "'0.06114"
$ws.Range("D6").Value  = "'3.558"
$ws.Range("D7").Value  = "'6.507"
$ws.Range("D8").Value  = "'1.338"
$ws.Range("D9").Value  = "'0.8237"
$ws.Range("D10").Value = "'0.01334"
$ws.Range("D12").Value = "'0.08094"
$ws.Range("D13").Value = "'0.03444"
$ws.Range("D14").Value = "'0.03183"
$ws.Range("D15").Value = "'0.09226"
$ws.Range("D16").Value = "'3.757"
$ws.Range("D17").Value = "'0.001650"
$ws.Range("D18").Value = "'0.04615"
$ws.Range("D19").Value = "'0.006344"
$ws.Range("D20").Value = "'0.006145"
$ws.Range("D21").Value = "'0.001069"
$ws.Range("D23").Value = "'3.731"
$ws.Range("D24").Value = "'2.290"
$ws.Range("D26").Value = "'0.1243"
$ws.Range("D28").Value = "'0.0002715"
$ws.Range("D40").Value = "'0.04598"
$ws.Range("D41").Value = "'0.007000"
$ws.Range("D42").Value = "'0.004003"
$ws.Range("D43").Value = "'0.1116"
$ws.Range("D44").Value = "'0.01154"
$ws.Range("D45").Value = "'0.00006053"
$ws.Range("D46").Value = "'0.0009907"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D48").Value = "'0.8031"
$ws.Range("D50").Value = "'0.00001902"
